$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.742.24"
$ws.Range("E2").Value = "  -2.18%  "
$ws.Range("D3").Value = "1.746.55"
$ws.Range("E3").Value = "  -4.56%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.004"
$ws.Range("E4").Value = "  +0.24%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "319.83"
$ws.Range("E5").Value = "  -3.03%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.003"
$ws.Range("E6").Value = "  +0.25%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4187"
$ws.Range("E7").Value = "  -5.94%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3643"
$ws.Range("E8").Value = "  -3.27%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "43.05"
$ws.Range("E9").Value = "  -3.80%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07400"
$ws.Range("E10").Value = "  -4.43%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.082"
$ws.Range("E11").Value = "  -4.47%  "
$ws.Range("E12").Value = "  +0.23%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "20.42"
$ws.Range("E13").Value = "  -8.46%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.033"
$ws.Range("E14").Value = "  -4.88%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.257"
$ws.Range("E15").Value = "  -4.18%  "
$ws.Range("D16").Value = "1.775.53"
$ws.Range("E16").Value = "  -3.00%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "90.63"
$ws.Range("E17").Value = "  -3.15%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.00001046"
$ws.Range("E18").Value = "  -3.55%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06324"
$ws.Range("E19").Value = "  -2.58%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "1.003"
$ws.Range("E20").Value = "  +0.24%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "16.94"
$ws.Range("E21").Value = "  -3.54%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.928"
$ws.Range("E22").Value = "  -6.67%  "
$ws.Range("D23").Value = "27.759.32"
$ws.Range("E23").Value = "  -2.26%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.16"
$ws.Range("E24").Value = "  -4.86%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.079"
$ws.Range("E25").Value = "  -6.60%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "156.95"
$ws.Range("E26").Value = "  +0.84%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "20.08"
$ws.Range("E27").Value = "  -3.40%  "
$ws.Range("D28").Value = "1.967.20"
$ws.Range("E28").Value = "  -3.82%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.125"
$ws.Range("E29").Value = "  -10.20%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "123.46"
$ws.Range("E30").Value = "  -4.13%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.126"
$ws.Range("E31").Value = "  -6.48%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.639"
$ws.Range("E32").Value = "  -0.99%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.542"
$ws.Range("E33").Value = "  -6.51%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.08798"
$ws.Range("E34").Value = "  -5.39%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "12.26"
$ws.Range("E35").Value = "  -7.19%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.02267"
$ws.Range("E36").Value = "  -3.73%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.2089"
$ws.Range("E37").Value = "  -4.85%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.05971"
$ws.Range("E38").Value = "  -4.22%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "4.935"
$ws.Range("E39").Value = "  -5.26%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.6277"
$ws.Range("E40").Value = "  -5.14%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.171"
$ws.Range("E41").Value = "  -2.79%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.002"
$ws.Range("E42").Value = "  +0.23%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.393"
$ws.Range("E43").Value = "  +0.05%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "7.739"
$ws.Range("E44").Value = "  -5.08%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "13.34"
$ws.Range("E45").Value = "  -4.58%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.5837"
$ws.Range("E46").Value = "  -4.63%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.667"
$ws.Range("E47").Value = "  -3.05%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "122.02"
$ws.Range("E48").Value = "  -4.31%  "
$ws.Range("E49").Value = "  -4.89%  "
$ws.Range("E50").Value = "  +0.95%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.06794"
$ws.Range("E51").Value = "  -3.13%  "
